$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F8").Value = 108
$ws1.Range("F18").Value = 111
$ws1.Range("F19").Value = 3955
$ws1.Range("F20").Value = 6284
$ws1.Range("F33").Value = 140
$ws1.Range("F35").Value = 304
$ws1.Range("F39").Value = 1561

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F9").Value = 108
$ws4.Range("F19").Value = 111
$ws4.Range("F20").Value = 3955
$ws4.Range("F21").Value = 6284
$ws4.Range("F23").Value = 38
$ws4.Range("F34").Value = 140
$ws4.Range("F36").Value = 304
$ws4.Range("F38").Value = 0
$ws4.Range("F40").Value = 1561
$ws4.Range("F44").Value = 57
